$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B17").Value = "Could I try projections on ice out dates?"
$ws.Range("B17").Select()
